$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 711.8182
$ws.Range("I12").Value = 661
$ws.Range("J12").Value = 800.75
$ws.Range("K12").Value = 661
$ws.Range("L12").Value = 800.75
$ws.Range("M12").Value = -491
$ws.Range("N12").Value = -1140.75
$ws.Range("H28").Value = 1514.1111
$ws.Range("I28").Value = 161
$ws.Range("K28").Value = 161
$ws.Range("M28").Value = 324
$ws.Range("H41").Value = 417.5
$ws.Range("I41").Value = 275
$ws.Range("K41").Value = 275
$ws.Range("M41").Value = 165
$ws.Range("H62").Value = 1801.75
$ws.Range("I62").Value = 1412.3334
$ws.Range("J62").Value = 2970
$ws.Range("K62").Value = 1412.3334
$ws.Range("L62").Value = 2970
$ws.Range("M62").Value = -788.3334
$ws.Range("N62").Value = -4218
$ws.Range("H65").Value = 1801.75
$ws.Range("I65").Value = 1412.3334
$ws.Range("J65").Value = 2970
$ws.Range("K65").Value = 7061.666999999999
$ws.Range("L65").Value = 14850
$ws.Range("M65").Value = -3941.666999999999
$ws.Range("N65").Value = -21090
$ws.Range("H86").Value = 1471.5714
$ws.Range("J86").Value = 1901.1666
$ws.Range("L86").Value = 1901.1666
$ws.Range("N86").Value = -4147.1666
$ws.Range("H89").Value = 1471.5714
$ws.Range("J89").Value = 1901.1666
$ws.Range("L89").Value = 9505.833000000001
$ws.Range("N89").Value = -20737.833
$ws.Range("H98").Value = 1653.625
$ws.Range("I98").Value = 1563.8667
$ws.Range("K98").Value = 1563.8667
$ws.Range("M98").Value = -65.86670000000004
$ws.Range("H106").Value = 4334
$ws.Range("I106").Value = 2856.2856
$ws.Range("K106").Value = 2856.2856
$ws.Range("M106").Value = -2225.2856
$ws.Range("H107").Value = 839.7857
$ws.Range("J107").Value = 2003
$ws.Range("L107").Value = 2003
$ws.Range("N107").Value = -5843
$ws.Range("H122").Value = 1653.625
$ws.Range("I122").Value = 1563.8667
$ws.Range("K122").Value = 4691.6001
$ws.Range("M122").Value = -2241.6001
$ws.Range("H138").Value = 2281.0527
$ws.Range("J138").Value = 2180
$ws.Range("L138").Value = 6540
$ws.Range("N138").Value = -16820

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3438.8235
$ws.Range("I32").Value = 2133.3684
$ws.Range("J32").Value = 7254.769
$ws.Range("K32").Value = 2133.3684
$ws.Range("L32").Value = 7254.769
$ws.Range("M32").Value = -1846.3684
$ws.Range("N32").Value = -7828.769
$ws.Range("H74").Value = 4518.615
$ws.Range("I74").Value = 4531
$ws.Range("K74").Value = 4531
$ws.Range("M74").Value = -3657
$ws.Range("H77").Value = 4518.615
$ws.Range("I77").Value = 4531
$ws.Range("K77").Value = 22655
$ws.Range("M77").Value = -18287
$ws.Range("H132").Value = 2178.1072
$ws.Range("I132").Value = 1294.7333
$ws.Range("J132").Value = 3197.3845
$ws.Range("K132").Value = 3884.199900000001
$ws.Range("L132").Value = 9592.1535
$ws.Range("M132").Value = -1354.199900000001
$ws.Range("N132").Value = -14652.1535

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 57041
$ws.Range("J36").Value = 57041
$ws.Range("L36").Value = 57041
$ws.Range("N36").Value = -58109
$ws.Range("H86").Value = 145121.36
$ws.Range("I86").Value = 2308.25
$ws.Range("J86").Value = 1002000
$ws.Range("K86").Value = 2308.25
$ws.Range("L86").Value = 1002000
$ws.Range("M86").Value = -1185.25
$ws.Range("N86").Value = -1004246
$ws.Range("H89").Value = 145121.36
$ws.Range("I89").Value = 2308.25
$ws.Range("J89").Value = 1002000
$ws.Range("K89").Value = 11541.25
$ws.Range("L89").Value = 5010000
$ws.Range("M89").Value = -5925.25
$ws.Range("N89").Value = -5021232
$ws.Range("H94").Value = 621.3333
$ws.Range("I94").Value = 488.26666
$ws.Range("J94").Value = 1286.6666
$ws.Range("K94").Value = 488.26666
$ws.Range("L94").Value = 1286.6666
$ws.Range("M94").Value = -37.26666
$ws.Range("N94").Value = -2188.6666
$ws.Range("H105").Value = 2419.8572
$ws.Range("I105").Value = 2419.8572
$ws.Range("K105").Value = 2419.8572
$ws.Range("M105").Value = -672.8571999999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H134").Value = 9126
$ws.Range("I134").Value = 10627.211
$ws.Range("J134").Value = 5051.2856
$ws.Range("K134").Value = 31881.633
$ws.Range("L134").Value = 15153.8568
$ws.Range("M134").Value = -29346.633
$ws.Range("N134").Value = -20223.8568

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1837.3334
$ws.Range("I86").Value = 1837.3334
$ws.Range("K86").Value = 1837.3334
$ws.Range("M86").Value = -714.3334
$ws.Range("H89").Value = 1837.3334
$ws.Range("I89").Value = 1837.3334
$ws.Range("K89").Value = 9186.666999999999
$ws.Range("M89").Value = -3570.666999999999
$ws.Range("H96").Value = 32750
$ws.Range("J96").Value = 32750
$ws.Range("L96").Value = 32750
$ws.Range("N96").Value = -38242
$ws.Range("H132").Value = 2012.8889
$ws.Range("I132").Value = 933.8570999999999
$ws.Range("J132").Value = 3523.5334
$ws.Range("K132").Value = 2801.5713
$ws.Range("L132").Value = 10570.6002
$ws.Range("M132").Value = -271.5712999999996
$ws.Range("N132").Value = -15630.6002
$ws.Range("H134").Value = 1063.1
$ws.Range("I134").Value = 1047.8889
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 3143.6667
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -608.6666999999998
$ws.Range("N134").Value = -8670

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13997.194
$ws.Range("I4").Value = 13997.194
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 41991.58199999999
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -41879.58199999999
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = -2831
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = -2898
$ws.Range("H68").Value = 2001.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2001.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6004.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7626.5
$ws.Range("H71").Value = 2001.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2001.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 18013.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -26125.5
$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -9564
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -26820
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 2140
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 280
$ws.Range("K86").Value = 12000
$ws.Range("L86").Value = 840
$ws.Range("M86").Value = -10814
$ws.Range("N86").Value = -3212
$ws.Range("H89").Value = 2140
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 280
$ws.Range("K89").Value = 36000
$ws.Range("L89").Value = 2520
$ws.Range("M89").Value = -30072
$ws.Range("N89").Value = -14376

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 55629.316
$ws.Range("J126").Value = 251996.25
$ws.Range("L126").Value = 755988.75
$ws.Range("N126").Value = -760928.75
$ws.Range("H132").Value = 3440.4814
$ws.Range("I132").Value = 2371.3333
$ws.Range("J132").Value = 4776.9165
$ws.Range("K132").Value = 7113.999899999999
$ws.Range("L132").Value = 14330.7495
$ws.Range("M132").Value = -4583.999899999999
$ws.Range("N132").Value = -19390.7495

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 3649
$ws.Range("I56").Value = 3973.5
$ws.Range("K56").Value = 3973.5
$ws.Range("M56").Value = -3282.5
$ws.Range("H122").Value = 4112.2188
$ws.Range("I122").Value = 1671.0625
$ws.Range("K122").Value = 5013.1875
$ws.Range("M122").Value = -2563.1875
$ws.Range("H132").Value = 2036.5
$ws.Range("I132").Value = 1697.2
$ws.Range("K132").Value = 5091.6
$ws.Range("M132").Value = -2561.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 48788.035
$ws.Range("I122").Value = 59453.91
$ws.Range("K122").Value = 178361.73
$ws.Range("M122").Value = -175911.73
$ws.Range("H126").Value = 4141.1113
$ws.Range("I126").Value = 3190.8096
$ws.Range("K126").Value = 9572.4288
$ws.Range("M126").Value = -7102.4288
